# Staff hours are now fortnightly.
# Staff sheet: "Average Weekly Hours" column is dropped, the single week of
# Mon-Fri start/end columns becomes two weeks (week 1 + week 2), and the
# "Hew Level" header is shortened to "Hew".

$wb = $excel.ActiveWorkbook
$staff = $wb.Worksheets.Item("Staff")

# --- Header row (row 1) ------------------------------------------------
$staff.Range("B1").Value = "Hew"

# Wipe out the old D1:R3 block (Average Weekly Hours + single week +
# leftover unlabeled/blank columns) so we can rebuild it cleanly as two
# full weeks (D:W).
$staff.Range("D1:R3").Clear()

$week1Headers = @("Mon 1 start","Mon 1 end","Tue 1 start","Tue 1 end","Wed 1 start","Wed 1 end","Thu 1 start","Thu 1 end","Fri 1 start","Fri 1 end")
$week2Headers = @("Mon 2 start","Mon 2 end","Tue 2 start","Tue 2 end","Wed 2 start","Wed 2 End","Thu 2 start","Thu 2 end","Fri 2 start","Fri 2 end")

$week1Cols = @("D","E","F","G","H","I","J","K","L","M")
$week2Cols = @("N","O","P","Q","R","S","T","U","V","W")

for ($i = 0; $i -lt 10; $i++) {
    $c1 = $staff.Range($week1Cols[$i] + "1")
    $c1.Value = $week1Headers[$i]
    $c1.NumberFormat = "HH:MM:SS"
    $c1.Font.Bold = $true

    $c2 = $staff.Range($week2Cols[$i] + "1")
    $c2.Value = $week2Headers[$i]
    $c2.NumberFormat = "HH:MM:SS"
    $c2.Font.Bold = $true
}

# --- Row 2: Edwina ------------------------------------------------------
$staff.Range("C2").Value = "y"

$row2 = @{
    "D" = 0.333333333333333; "E" = 0.708333333333333
    "F" = 0.333333333333333; "G" = 0.708333333333333
    "H" = 0.354166666666667; "I" = 0.708333333333333
    "J" = 0.333333333333333; "K" = 0.6875
    "N" = 0.333333333333333; "O" = 0.708333333333333
    "P" = 0.333333333333333; "Q" = 0.708333333333333
    "R" = 0.354166666666667; "S" = 0.708333333333333
    "T" = 0.333333333333333; "U" = 0.6875
    "V" = 0.333333333333333; "W" = 0.6875
}
foreach ($col in $row2.Keys) {
    $cell = $staff.Range($col + "2")
    $cell.Value = $row2[$col]
    $cell.NumberFormat = "HH:MM:SS"
}
# Fri 1 (L2/M2) stays empty for Edwina - leave untouched (no cell at all).

# --- Row 3: Rowena -------------------------------------------------------
$staff.Range("C3").Value = "n"

# Blank-but-time-formatted placeholder cells (no value).
$blankCols3 = @("D","E","F","G","H","I","N","O","P","Q")
foreach ($col in $blankCols3) {
    $cell = $staff.Range($col + "3")
    $cell.ClearContents()
    $cell.NumberFormat = "HH:MM:SS"
}

$row3 = @{
    "J" = 0.479166666666667; "K" = 0.833333333333333
    "L" = 0.354166666666667; "M" = 0.708333333333333
    "R" = 0.479166666666667; "S" = 0.833333333333333
    "T" = 0.354166666666667; "U" = 0.708333333333333
}
foreach ($col in $row3.Keys) {
    $cell = $staff.Range($col + "3")
    $cell.Value = $row3[$col]
    $cell.NumberFormat = "HH:MM:SS"
}
# Fri 2 (V3/W3) stays empty for Rowena - leave untouched (no cell at all).

# Restore the selection on the Staff sheet.
$staff.Range("G14").Select()
